$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.104028125842603
$ws.Range("C2").Value = 0.08797166765044295
$ws.Range("D2").Value = 0.08489952917120291
$ws.Range("F2").Value = 2.656232113367636
$ws.Range("G2").Value = 1.935723986379486
$ws.Range("H2").Value = 1.671791366903335
$ws.Range("J2").Value = 0.2275797636615149
$ws.Range("K2").Value = 0.6070571432156555
$ws.Range("L2").Value = 0.304746940193894
$ws.Range("M2").Value = 0.2965859752951197
$ws.Range("N2").Value = 3.09847412006097

$ws.Range("B3").Value = 1.071856801390254
$ws.Range("C3").Value = 0.08576388632950938
$ws.Range("D3").Value = 0.0842236403907819
$ws.Range("F3").Value = 2.660963749101271
$ws.Range("G3").Value = 1.938358767601315
$ws.Range("H3").Value = 1.677969308742206
$ws.Range("J3").Value = 0.2286467372366587
$ws.Range("K3").Value = 0.57526942377811
$ws.Range("L3").Value = 0.3032378532016509
$ws.Range("M3").Value = 0.2911560665213138
$ws.Range("N3").Value = 3.12103516581643

$ws.Range("B4").Value = 1.05263721539086
$ws.Range("C4").Value = 0.08438641722565166
$ws.Range("D4").Value = 0.08383505374934686
$ws.Range("F4").Value = 2.664911076751785
$ws.Range("G4").Value = 1.940818397660834
$ws.Range("H4").Value = 1.682329982602269
$ws.Range("J4").Value = 0.2293488229790892
$ws.Range("K4").Value = 0.5560448379864056
$ws.Range("L4").Value = 0.302422742825307
$ws.Range("M4").Value = 0.2879562679319534
$ws.Range("N4").Value = 3.135672147332535

$ws.Range("B5").Value = 1.04493982934693
$ws.Range("C5").Value = 0.08381957732726164
$ws.Range("D5").Value = 0.08368338013033139
$ws.Range("F5").Value = 2.666781959650379
$ws.Range("G5").Value = 1.942032517531871
$ws.Range("H5").Value = 1.684249835162603
$ws.Range("J5").Value = 0.2296467589935478
$ws.Range("K5").Value = 0.5482847504474933
$ws.Range("L5").Value = 0.3021186987232198
$ws.Range("M5").Value = 0.2866861897734552
$ws.Range("N5").Value = 3.141834151025989

$ws.Range("B6").Value = 1.04366983879558
$ws.Range("C6").Value = 0.08372512098255669
$ws.Range("D6").Value = 0.08365859930566444
$ws.Range("F6").Value = 2.667108469607008
$ws.Range("G6").Value = 1.942246916642489
$ws.Range("H6").Value = 1.68457725749046
$ws.Range("J6").Value = 0.2296969462148866
$ws.Range("K6").Value = 0.5470006791163087
$ws.Range("L6").Value = 0.3020699134223932
$ws.Range("M6").Value = 0.2864773433042025
$ws.Range("N6").Value = 3.142869266453644

$ws.Range("B7").Value = 1.052532859432517
$ws.Range("C7").Value = 0.08437879494362477
$ws.Range("D7").Value = 0.08383298113434279
$ws.Range("F7").Value = 2.664935245668978
$ws.Range("G7").Value = 1.940833913993586
$ws.Range("H7").Value = 1.682355295822347
$ws.Range("J7").Value = 0.2293527931183377
$ws.Range("K7").Value = 0.5559398822790627
$ws.Range("L7").Value = 0.3024185284050986
$ws.Range("M7").Value = 0.2879390019288088
$ws.Range("N7").Value = 3.135754451261676

$ws.Range("B8").Value = 1.092824986308813
$ws.Range("C8").Value = 0.08721496128974593
$ws.Range("D8").Value = 0.08466102283201593
$ws.Range("F8").Value = 2.657647456774939
$ws.Range("G8").Value = 1.936457790422679
$ws.Range("H8").Value = 1.673803857619333
$ws.Range("J8").Value = 0.2279379228964231
$ws.Range("K8").Value = 0.5960360875088782
$ws.Range("L8").Value = 0.3042035223971737
$ws.Range("M8").Value = 0.2946859636636034
$ws.Range("N8").Value = 3.106090407847638

$ws.Range("B9").Value = 1.176054239058573
$ws.Range("C9").Value = 0.09260360142836532
$ws.Range("D9").Value = 0.08649289005057881
$ws.Range("F9").Value = 2.651614531716291
$ws.Range("G9").Value = 1.934553113312418
$ws.Range("H9").Value = 1.661529825578214
$ws.Range("J9").Value = 0.2255349872779515
$ws.Range("K9").Value = 0.6769798817778394
$ws.Range("L9").Value = 0.3085849484840111
$ws.Range("M9").Value = 0.3089770210457203
$ws.Range("N9").Value = 3.054139473129595

$ws.Range("B10").Value = 1.239756174500172
$ws.Range("C10").Value = 0.09645819980519832
$ws.Range("D10").Value = 0.08796383074128755
$ws.Range("F10").Value = 2.652206103544657
$ws.Range("G10").Value = 1.937222876470528
$ws.Range("H10").Value = 1.655244517861306
$ws.Range("J10").Value = 0.2239947529891673
$ws.Range("K10").Value = 0.7378525181954672
$ws.Range("L10").Value = 0.312337127240724
$ws.Range("M10").Value = 0.3201185374921991
$ws.Range("N10").Value = 3.019756968991146

$ws.Range("B11").Value = 1.26928714593015
$ws.Range("C11").Value = 0.09818932645012524
$ws.Range("D11").Value = 0.08865979414333935
$ws.Range("F11").Value = 2.653564125565623
$ws.Range("G11").Value = 1.939320934677895
$ws.Range("H11").Value = 1.652976880857878
$ws.Range("J11").Value = 0.223342677136575
$ws.Range("K11").Value = 0.7658484367845801
$ws.Range("L11").Value = 0.3141590614596907
$ws.Range("M11").Value = 0.3253256187073958
$ws.Range("N11").Value = 3.004936255827197

$ws.Range("B12").Value = 1.280548794299506
$ws.Range("C12").Value = 0.09884165754232299
$ws.Range("D12").Value = 0.088927160422692
$ws.Range("F12").Value = 2.654234744997211
$ws.Range("G12").Value = 1.940242422513222
$ws.Range("H12").Value = 1.652203112218046
$ws.Range("J12").Value = 0.223102717597687
$ws.Range("K12").Value = 0.7764933066863193
$ws.Range("L12").Value = 0.3148654452451751
$ws.Range("M12").Value = 0.327317253979011
$ws.Range("N12").Value = 2.99944190429526

$ws.Range("B13").Value = 1.278119896156767
$ws.Range("C13").Value = 0.09870130901416019
$ws.Range("D13").Value = 0.08886940900719509
$ws.Range("F13").Value = 2.654083364145492
$ws.Range("G13").Value = 1.94003831622183
$ws.Range("H13").Value = 1.652365981667373
$ws.Range("J13").Value = 0.2231540875994185
$ws.Range("K13").Value = 0.7741988183717297
$ws.Range("L13").Value = 0.3147125823206665
$ws.Range("M13").Value = 0.3268874398884023
$ws.Range("N13").Value = 3.000619966043431

$ws.Range("B14").Value = 1.270212070236511
$ws.Range("C14").Value = 0.09824305840729153
$ws.Range("D14").Value = 0.08868171418201598
$ws.Range("F14").Value = 2.653616165268701
$ws.Range("G14").Value = 1.939394201118773
$ws.Range("H14").Value = 1.652911520949942
$ws.Range("J14").Value = 0.2233227959914821
$ws.Range("K14").Value = 0.7667233285809516
$ws.Range("L14").Value = 0.3142168467576454
$ws.Range("M14").Value = 0.3254890751022543
$ws.Range("N14").Value = 3.004481868656121

$ws.Range("B15").Value = 1.265378558414426
$ws.Range("C15").Value = 0.09796194891530519
$ws.Range("D15").Value = 0.08856724205818267
$ws.Range("F15").Value = 2.653350349475446
$ws.Range("G15").Value = 1.939016199313656
$ws.Range("H15").Value = 1.653256736818079
$ws.Range("J15").Value = 0.2234270415583062
$ws.Range("K15").Value = 0.7621500187706545
$ws.Range("L15").Value = 0.3139153349961674
$ws.Range("M15").Value = 0.3246351145876432
$ws.Range("N15").Value = 3.006862753871637

$ws.Range("B16").Value = 1.237837328848059
$ws.Range("C16").Value = 0.09634461729285704
$ws.Range("D16").Value = 0.08791888448703844
$ws.Range("F16").Value = 2.65213924661964
$ws.Range("G16").Value = 1.937103538171485
$ws.Range("H16").Value = 1.655404606049302
$ws.Range("J16").Value = 0.2240383436996751
$ws.Range("K16").Value = 0.7360290164932621
$ws.Range("L16").Value = 0.3122203671886581
$ws.Range("M16").Value = 0.3197810237606333
$ws.Range("N16").Value = 3.020742039700039

$ws.Range("B17").Value = 1.221082823555236
$ws.Range("C17").Value = 0.09534671594254718
$ws.Range("D17").Value = 0.08752798202962708
$ws.Range("F17").Value = 2.651675015100082
$ws.Range("G17").Value = 1.936156444632516
$ws.Range("H17").Value = 1.656873677742254
$ws.Range("J17").Value = 0.2244257878730664
$ws.Range("K17").Value = 0.7200823819132722
$ws.Range("L17").Value = 0.3112099641613852
$ws.Range("M17").Value = 0.3168386440600273
$ws.Range("N17").Value = 3.02946657318904

$ws.Range("B18").Value = 1.211498125852131
$ws.Range("C18").Value = 0.09477064555614589
$ws.Range("D18").Value = 0.08730567240503007
$ws.Range("F18").Value = 2.651510507998353
$ws.Range("G18").Value = 1.935694870586758
$ws.Range("H18").Value = 1.657774339141497
$ws.Range("J18").Value = 0.224653209634857
$ws.Range("K18").Value = 0.7109389977000262
$ws.Range("L18").Value = 0.3106396412406127
$ws.Range("M18").Value = 0.3151593323222244
$ws.Range("N18").Value = 3.034561881391767

$ws.Range("B19").Value = 1.208261868673816
$ws.Range("C19").Value = 0.09457523652164923
$ws.Range("D19").Value = 0.08723083736524728
$ws.Range("F19").Value = 2.651472421012755
$ws.Range("G19").Value = 1.935552875694398
$ws.Range("H19").Value = 1.658088856524046
$ws.Range("J19").Value = 0.2247309969258797
$ws.Range("K19").Value = 0.7078481455987742
$ws.Range("L19").Value = 0.3104484031373289
$ws.Range("M19").Value = 0.3145929943578167
$ws.Range("N19").Value = 3.036300322037654

$ws.Range("B20").Value = 1.222860986253949
$ws.Range("C20").Value = 0.09545316195777787
$ws.Range("D20").Value = 0.08756933298645464
$ws.Range("F20").Value = 2.6517138258736
$ws.Range("G20").Value = 1.936248657082615
$ws.Range("H20").Value = 1.656711529784033
$ws.Range("J20").Value = 0.2243840705127162
$ws.Range("K20").Value = 0.7217769610208222
$ws.Range("L20").Value = 0.3113164026241719
$ws.Range("M20").Value = 0.3171505139006996
$ws.Range("N20").Value = 3.028529843062003

$ws.Range("B21").Value = 1.272532653243132
$ws.Range("C21").Value = 0.09837774477189498
$ws.Range("D21").Value = 0.08873674132411224
$ws.Range("F21").Value = 2.653749150909022
$ws.Range("G21").Value = 1.939579947082137
$ws.Range("H21").Value = 1.652748978739922
$ws.Range("J21").Value = 0.2232730533032736
$ws.Range("K21").Value = 0.7689178857205548
$ws.Range("L21").Value = 0.3143620103684981
$ws.Range("M21").Value = 0.3258992716462146
$ws.Range("N21").Value = 3.003344333559426

$ws.Range("B22").Value = 1.305455650285296
$ws.Range("C22").Value = 0.1002704264071781
$ws.Range("D22").Value = 0.08952195906235971
$ws.Range("F22").Value = 2.655990673763782
$ws.Range("G22").Value = 1.942497357503512
$ws.Range("H22").Value = 1.650654235756406
$ws.Range("J22").Value = 0.2225875421906025
$ws.Range("K22").Value = 0.7999801295231066
$ws.Range("L22").Value = 0.316448373888278
$ws.Range("M22").Value = 0.3317325972734935
$ws.Range("N22").Value = 2.987571544453004

$ws.Range("B23").Value = 1.287842149857511
$ws.Range("C23").Value = 0.09926197627785882
$ws.Range("D23").Value = 0.08910085010925428
$ws.Range("F23").Value = 2.654711016719858
$ws.Range("G23").Value = 1.940872569222762
$ws.Range("H23").Value = 1.651726989149282
$ws.Range("J23").Value = 0.2229497034281245
$ws.Range("K23").Value = 0.7833786101245437
$ws.Range("L23").Value = 0.3153260978460821
$ws.Range("M23").Value = 0.3286087141093219
$ws.Range("N23").Value = 2.995926878682162

$ws.Range("B24").Value = 1.222056930304518
$ws.Range("C24").Value = 0.09540504507767622
$ws.Range("D24").Value = 0.08755063066113422
$ws.Range("F24").Value = 2.651695960575751
$ws.Range("G24").Value = 1.93620670952653
$ws.Range("H24").Value = 1.65678466219812
$ws.Range("J24").Value = 0.2244029163676835
$ws.Range("K24").Value = 0.7210107653123998
$ws.Range("L24").Value = 0.3112682488609693
$ws.Range("M24").Value = 0.3170094792184344
$ws.Range("N24").Value = 3.028953091235813

$ws.Range("B25").Value = 1.153088884474727
$ws.Range("C25").Value = 0.091164248905514
$ws.Range("D25").Value = 0.08597522998066154
$ws.Range("F25").Value = 2.652363579527787
$ws.Range("G25").Value = 1.934353709341011
$ws.Range("H25").Value = 1.664369858909623
$ws.Range("J25").Value = 0.2261454003727295
$ws.Range("K25").Value = 0.6548354262845635
$ws.Range("L25").Value = 0.3073057448710941
$ws.Range("M25").Value = 0.3049978668053868
$ws.Range("N25").Value = 3.067528313435943
